$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.540.69"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.515.61"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.561"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("D14").Value = "2.898.39"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "2.533.91"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.805"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("D18").Value = "42.484.29"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.44%  "
$ws.Range("D20").Value = "0.0₃0935"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.91%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.55%  "
$ws.Range("E28").Value = "  -4.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.16%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0781"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.107"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0294"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").Value = "2.000.81"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "2.746.88"
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "78.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.187"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.80%  "
